$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking text (e.g. "1.00", "0.330") need an
# explicit Text number format first, otherwise Excel COM auto-converts the
# assigned string into a real number (losing the original text representation,
# e.g. trailing zeros / exact decimal form), same as interactive Excel behavior.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '68.546.14'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '2.459.04'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '558.15'
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").Value = '160.65'
$ws.Range("E6").Value = '  -1.70%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '0.506'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.150'
$ws.Range("E9").Value = '  -0.96%  '
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = '4.85'
$ws.Range("E11").Value = '  +0.72%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = '0.330'
$ws.Range("E12").Value = '  -3.09%  '
$ws.Range("D13").Value = '68.453.86'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").Value = '0.0000167'
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").Value = '23.30'
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("D16").Value = '10.60'
$ws.Range("E16").Value = '  -3.31%  '
$ws.Range("D17").Value = '333.59'
$ws.Range("E17").Value = '  -3.00%  '
$ws.Range("D18").Value = '6.89'
$ws.Range("E18").Value = '  -3.64%  '
$ws.Range("D19").Value = '3.77'
$ws.Range("E19").Value = '  -1.27%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '1.00'
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("B21").Value = 'SuiNetwork'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D21").Value = '1.87'
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("D22").Value = '66.37'
$ws.Range("E22").Value = '  -2.51%  '
$ws.Range("D23").Value = '3.62'
$ws.Range("E23").Value = '  -3.19%  '
$ws.Range("D24").Value = '8.12'
$ws.Range("E24").Value = '  -1.27%  '
$ws.Range("D25").Value = '0.0₃0812'
$ws.Range("E25").Value = '  -3.32%  '
$ws.Range("D26").Value = '7.17'
$ws.Range("E26").Value = '  -1.79%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").Value = '426.60'
$ws.Range("E28").Value = '  -2.10%  '
$ws.Range("D29").Value = '1.13'
$ws.Range("E29").Value = '  -3.90%  '
$ws.Range("E30").Value = '  -4.36%  '
$ws.Range("D31").Value = '157.56'
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("D32").Value = '18.99'
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").Value = '0.109'
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("D35").Value = '17.70'
$ws.Range("E35").Value = '  -1.11%  '
$ws.Range("E36").Value = '  -2.35%  '
$ws.Range("D37").Value = '4.38'
$ws.Range("E37").Value = '  -2.38%  '
$ws.Range("D38").Value = '1.45'
$ws.Range("E38").Value = '  -4.90%  '
$ws.Range("D39").Value = '1.07'
$ws.Range("E39").Value = '  -3.05%  '
$ws.Range("D40").Value = '2.04'
$ws.Range("E40").Value = '  -1.99%  '
$ws.Range("D41").Value = '3.33'
$ws.Range("E41").Value = '  -1.27%  '
$ws.Range("D42").Value = '128.92'
$ws.Range("E42").Value = '  -4.36%  '
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("D44").Value = '0.480'
$ws.Range("E44").Value = '  -1.27%  '
$ws.Range("D45").Value = '0.559'
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("D46").Value = '0.0908'
$ws.Range("E46").Value = '  -0.76%  '
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("D48").Value = '1.38'
$ws.Range("E48").Value = '  -3.74%  '
$ws.Range("D49").Value = '4.90'
$ws.Range("E49").Value = '  -8.88%  '
$ws.Range("D50").Value = '16.71'
$ws.Range("E50").Value = '  -5.30%  '
$ws.Range("D51").Value = '0.0₆0205'
$ws.Range("E51").Value = '  +0.12%  '
